$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4908
$ws1.Range("F5").Value = 799
$ws1.Range("F6").Value = 246
$ws1.Range("F7").Value = 1275
$ws1.Range("F10").Value = 214
$ws1.Range("F11").Value = 91
$ws1.Range("F13").Value = 164
$ws1.Range("F15").Value = 4292
$ws1.Range("F16").Value = 6591
$ws1.Range("F18").Value = 47
$ws1.Range("F19").Value = 91
$ws1.Range("F20").Value = 554
$ws1.Range("F22").Value = 4062
$ws1.Range("F23").Value = 426
$ws1.Range("F24").Value = 58
$ws1.Range("F25").Value = 35
$ws1.Range("F26").Value = 2648
$ws1.Range("F27").Value = 571
$ws1.Range("F30").Value = 327
$ws1.Range("F31").Value = 337
$ws1.Range("F32").Value = 388
$ws1.Range("F33").Value = 202
$ws1.Range("F34").Value = 23
$ws1.Range("F35").Value = 1593
$ws1.Range("F36").Value = 1000
$ws1.Range("F38").Value = 107
$ws1.Range("F39").Value = 70
$ws1.Range("F40").Value = 516
$ws1.Range("F41").Value = 492
$ws1.Range("F44").Value = 608

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4908
$ws4.Range("F6").Value = 799
$ws4.Range("F7").Value = 246
$ws4.Range("F8").Value = 1275
$ws4.Range("F11").Value = 214
$ws4.Range("F12").Value = 91
$ws4.Range("F14").Value = 164
$ws4.Range("F16").Value = 4292
$ws4.Range("F17").Value = 6591
$ws4.Range("F19").Value = 47
$ws4.Range("F20").Value = 91
$ws4.Range("F21").Value = 554
$ws4.Range("F23").Value = 4062
$ws4.Range("F24").Value = 426
$ws4.Range("F25").Value = 58
$ws4.Range("F26").Value = 35
$ws4.Range("F27").Value = 2648
$ws4.Range("F28").Value = 571
$ws4.Range("F31").Value = 327
$ws4.Range("F32").Value = 337
$ws4.Range("F33").Value = 388
$ws4.Range("F34").Value = 202
$ws4.Range("F35").Value = 23
$ws4.Range("F36").Value = 1593
$ws4.Range("F37").Value = 1000
$ws4.Range("F39").Value = 107
$ws4.Range("F40").Value = 70
$ws4.Range("F41").Value = 516
$ws4.Range("F42").Value = 492
$ws4.Range("F45").Value = 608
